# Update visitor/ticket counts on both "展览" and "全部类型" sheets.
# Both sheets contain the same event rows (2-17), and both need the
# identical set of F/G column updates below.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G3").Value = 50
    $ws.Range("F4").Value = 14
    $ws.Range("F6").Value = 67
    $ws.Range("F9").Value = 474
    $ws.Range("F10").Value = 6185
    $ws.Range("F11").Value = 158
    $ws.Range("F12").Value = 108
    $ws.Range("F13").Value = 999
    $ws.Range("F14").Value = 234
    $ws.Range("F15").Value = 78
    $ws.Range("F16").Value = 169
    $ws.Range("F17").Value = 388
}
